$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A3").Value = -22.1767
$ws.Range("A14").Value = -21.82899999999999
$ws.Range("A21").Value = -20.10409999999998
$ws.Range("A23").Value = -20.06969999999998
$ws.Range("A25").Value = -21.77499999999999
$ws.Range("A26").Value = -21.08419999999996
$ws.Range("A29").Value = -20.68129999999997
$ws.Range("A53").Value = -22.0549
$ws.Range("A57").Value = -22.58110000000001
$ws.Range("A59").Value = -22.2971
$ws.Range("A69").Value = -21.62179999999999
$ws.Range("A79").Value = -20.18830000000001
$ws.Range("A83").Value = -21.9583
$ws.Range("A91").Value = -20.49639999999998
$ws.Range("A93").Value = -21.42280000000001
$ws.Range("A103").Value = -21.76569999999999
